# Deploy the implementation guide.
#
# Updates the generated CodeSystem metadata sheet to the newly-published
# status/date, and turns on the (previously-defined-but-unapplied) wrap
# text / top vertical alignment formatting for the data grids on both
# worksheets.

$wb = $excel.ActiveWorkbook

$wsMeta     = $wb.Worksheets.Item("Metadata")
$wsConcepts = $wb.Worksheets.Item("Concepts")

# --- Metadata content: Status -> draft, Date -> new publication date ---
$wsMeta.Range("B6").Value = "draft"
$wsMeta.Range("B8").Value = "2023-08-01T16:12:28+00:00"

# --- Formatting: apply the vertical-top / wrap-text alignment that was ---
# --- already defined on these cell styles but not switched on          ---
$wsMeta.UsedRange.WrapText = $true
$wsConcepts.UsedRange.WrapText = $true
